$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.588.97'
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("D3").Value = '1.753.41'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3597'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07489'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.094'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.003'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("E15").Value = '  -3.67%  '
$ws.Range("D16").Value = '1.753.11'
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001063'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06409'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.836'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.63%  '
$ws.Range("D23").Value = '27.644.87'
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.107'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("D28").Value = '1.955.33'
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '127.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.082'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.077'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09210'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.663'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.522'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.91'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02291'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2104'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06031'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6343'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.955'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.04%  '
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.383'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.750'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5891'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.713'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.951'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.149'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.30%  '
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.87%  '
